$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last refreshed" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 20:22"

# --- Row 4: Estados Unidos (updated stats) ---
$ws.Range("B4").Value = 903775
$ws.Range("C4").Value = 17333
$ws.Range("E4").Value = 762526
$ws.Range("G4").Value = 752
$ws.Range("H4").Value = 50988

# --- Row 18: Suiza (updated stats) ---
$ws.Range("D18").Value = 21000
$ws.Range("E18").Value = 6099

# --- Row 22: Peru (updated stats) ---
$ws.Range("B22").Value = 21648
$ws.Range("C22").Value = 734
$ws.Range("D22").Value = 7496
$ws.Range("E22").Value = 13518
$ws.Range("F22").Value = 505
$ws.Range("G22").Value = 62
$ws.Range("H22").Value = 634

# --- Rows 55/56: Marruecos overtakes Luxemburgo in the ranking ---
# Row 55 becomes Marruecos with fresh figures; row 56 becomes Luxemburgo
# carrying the figures Marruecos's old slot (row 55) used to hold.
$ws.Range("A55").Value = "Marruecos"
$ws.Range("B55").Value = 3758
$ws.Range("C55").Value = 190
$ws.Range("D55").Value = 486
$ws.Range("E55").Value = 3114
$ws.Range("F55").Value = 1
$ws.Range("G55").Value = 3
$ws.Range("H55").Value = 158

$ws.Range("A56").Value = "Luxemburgo"
$ws.Range("B56").Value = 3695
$ws.Range("C56").Value = 30
$ws.Range("D56").Value = 3007
$ws.Range("E56").Value = 603
$ws.Range("F56").Value = 25
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = 85

# --- Rows 108/109: Georgia overtakes Jordania in the ranking ---
$ws.Range("A108").Value = "Georgia"
$ws.Range("B108").Value = 444
$ws.Range("C108").Value = 19
$ws.Range("D108").Value = 132
$ws.Range("E108").Value = 307
$ws.Range("F108").Value = 6
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 5

$ws.Range("A109").Value = "Jordania"
$ws.Range("B109").Value = 441
$ws.Range("C109").Value = 4
$ws.Range("D109").Value = 326
$ws.Range("E109").Value = 108
$ws.Range("F109").Value = 5
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 7

# --- Row 111: Sri Lanka (updated stats) ---
$ws.Range("B111").Value = 417
$ws.Range("C111").Value = 49
$ws.Range("E111").Value = 301
